$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A2").Value = 3013
$ws.Range("E2").Value = 46200608013
$ws.Range("AM2").Value = 1
$ws.Range("AN2").Value = 0
$ws.Range("AQ2").ClearContents()
$ws.Range("AU2").Value = "Thẻ BHYT hợp lệ"
$ws.Range("BA2").Value = "2024-05-09T09:13:34.2378979+07:00"
$ws.Range("BF2").Value = "||149|09/05/2024 09:13||||||3839||||New|4803|80|||2|1083660|||"

$ws.Range("AQ8").Select()
